$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds binary labels (0/1) for rows 2-41.
# The edit flips every label: rows 2-21 go 0 -> 1, rows 22-41 go 1 -> 0.
$ws.Range("B2:B21").Value = 1
$ws.Range("B22:B41").Value = 0

# Move the active selection from E10 to D10.
[void]$ws.Range("D10").Select()
